$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force text-typed values onto cells whose new text
# would otherwise be auto-converted to a number by plain Value assignment.
$ws.Range("Z1").NumberFormat = "@"

$ws.Range("D2").Value = "29.510.50"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.878.34"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("Z1").Value = "0.7188"
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").Value = "241.83"
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("Z1").Value = "0.07902"
$ws.Range("Z1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("Z1").Value = "0.3099"
$ws.Range("Z1").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("Z1").Value = "25.47"
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("Z1").Value = "0.08271"
$ws.Range("Z1").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "1.906.93"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("Z1").Value = "0.7289"
$ws.Range("Z1").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("Z1").Value = "5.282"
$ws.Range("Z1").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("Z1").Value = "91.30"
$ws.Range("Z1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "29.525.05"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("Z1").Value = "5.903"
$ws.Range("Z1").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("Z1").Value = "245.86"
$ws.Range("Z1").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +3.55%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("Z1").Value = "13.32"
$ws.Range("Z1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "2.128.22"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("Z1").Value = "8.074"
$ws.Range("Z1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +6.72%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +14.91%  "
$ws.Range("Z1").Value = "163.47"
$ws.Range("Z1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("Z1").Value = "18.32"
$ws.Range("Z1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("Z1").Value = "4.394"
$ws.Range("Z1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("Z1").Value = "4.112"
$ws.Range("Z1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("Z1").Value = "0.05214"
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("Z1").Value = "1.948"
$ws.Range("Z1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("Z1").Value = "0.7278"
$ws.Range("Z1").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("Z1").Value = "2.677"
$ws.Range("Z1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("Z1").Value = "0.01873"
$ws.Range("Z1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").Value = "1.206.37"
$ws.Range("E39").Value = "  +5.49%  "
$ws.Range("Z1").Value = "2.706"
$ws.Range("Z1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("Z1").Value = "0.9088"
$ws.Range("Z1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("Z1").Value = "6.180"
$ws.Range("Z1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +3.57%  "
$ws.Range("Z1").Value = "73.50"
$ws.Range("Z1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +4.34%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("Z1").Value = "102.46"
$ws.Range("Z1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "2.023.79"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("Z1").Value = "0.5295"
$ws.Range("Z1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("Z1").Value = "1.801"
$ws.Range("Z1").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("Z1").Value = "2.933"
$ws.Range("Z1").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +9.15%  "
$ws.Range("Z1").Value = "9.319"
$ws.Range("Z1").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("E51").Value = "  +1.94%  "

# Clean up the scratch cell so it leaves no trace in the used range.
$ws.Range("Z1").Clear()
